$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.155.77'
$ws.Range("E2").Value = '  +0.24%  '
$ws.Range("D3").Value = '2.515.20'
$ws.Range("E3").Value = '  +2.06%  '
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E4").Value = '  +0.06%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '520.61'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.37%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '132.35'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -0.57%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +0.21%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.557'
$c.Style = "Normal"
$ws.Range("E8").Value = '  +0.25%  '
$ws.Range("D9").Value = '2.513.37'
$ws.Range("E9").Value = '  +1.51%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.0972'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -0.66%  '
$ws.Range("E11").Value = '  -1.11%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '5.15'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -2.67%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.331'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -2.35%  '
$ws.Range("D14").Value = '2.958.12'
$ws.Range("E14").Value = '  +1.84%  '
$ws.Range("D15").Value = '58.321.93'
$ws.Range("E15").Value = '  +0.60%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '22.03'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -1.01%  '
$ws.Range("E17").Value = '  -0.27%  '
$ws.Range("D18").Value = '2.509.27'
$ws.Range("E18").Value = '  +1.43%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '10.59'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -0.36%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '320.75'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +0.03%  '
$ws.Range("E21").Value = '  -0.23%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '6.15'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +7.48%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -0.08%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '64.57'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -0.11%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '0.404'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -0.91%  '
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("E27").Value = '  -0.11%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '7.36'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +0.48%  '
$ws.Range("D29").Value = '0.0₃0752'
$ws.Range("E29").Value = '  +1.07%  '
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '1.71'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +1.43%  '
$ws.Range("B31").Value = 'Monero'
$ws.Range("C31").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '167.55'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -0.49%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '1.19'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +1.62%  '
$ws.Range("E33").Value = '  -0.18%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.997'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -0.12%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '18.03'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +0.04%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '1.25'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -7.61%  '
$ws.Range("E38").Value = '  -1.14%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '1.47'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +0.56%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '36.15'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -0.56%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.769'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -3.14%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '277.88'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +1.84%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '3.47'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +0.65%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '5.01'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +0.24%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.595'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +0.76%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '125.23'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +0.72%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.0919'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +1.46%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.0500'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +2.65%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '17.61'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +0.35%  '
$ws.Range("E50").Value = '  +0.23%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '16.76'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -0.48%  '
